# multiple_session_batch.xlsx:
#  - "fixes for new input_samples folder": the File-path value cells on the
#    "Sessions" sheet pointed at a "sample_inputs/" folder; the repo folder
#    was renamed to "input_samples/", so every such path cell is updated to
#    keep the same filename under the new folder name.
#  - the active selection / frozen-pane scroll position left over from the
#    editing session also moved (pane scrolled up to A10, selection moved to
#    C21:C23) and the workbook window was resized.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sessions")

# ---- File path cells: sample_inputs/... -> input_samples/... ----
$pathCells = @(
    "C12","D12",   # Manufacturers File
    "C13","D13",   # Market Classes File
    "C14","D14",   # Vehicles File
    "C15","D15",   # Demanded Shares File
    "C16","D16",   # Fuels File
    "C17","D17",   # Fuel Scenario Annual Data File
    "C19","D19",   # Cost File
    "C21",         # GHG Standards File (Footprint)
    "D21",         # GHG Standards File (Flat)
    "C22","D22",   # GHG Standards Fuels File
    "C23","D23"    # ZEV Requirement File
)

foreach ($addr in $pathCells) {
    $cell = $ws.Range($addr)
    $cell.Value = $cell.Value2 -replace '^sample_inputs/', 'input_samples/'
}

# ---- Window / selection state ----
$win = $excel.ActiveWindow

# Move the frozen-pane scroll so row 10 is the first visible row under the
# freeze at row 9, and select C21:C23 (active cell C21).
[void]$ws.Range("C21:C23").Select()
$win.ScrollRow = 10

# Resize / reposition the workbook window.
$win.Left = 255
$win.Top = 3255
$win.Width = 19815
$win.Height = 15180

[void]$wb.Save()
